$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the SmartArt shape, which may live inside a (possibly nested) group.
function Find-SmartArtShape($shapeRange) {
    for ($i = 1; $i -le $shapeRange.Count; $i++) {
        $sh = $shapeRange.Item($i)
        if ($sh.HasSmartArt) {
            return $sh
        }
        if ($sh.Type -eq 6) {
            $found = Find-SmartArtShape($sh.GroupItems)
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

$dia = Find-SmartArtShape($s.Shapes)
$sa = $dia.SmartArt
$nodes = $sa.AllNodes

# Two pairs of node labels were swapped with each other:
#   "חיפוש עמדה"   <-> "רישום שחקן לעמדה"
#   "מחיקת טכנאי"  <-> "עדכון טכנאי"
$pairs = @(
    @("חיפוש עמדה", "רישום שחקן לעמדה"),
    @("מחיקת טכנאי", "עדכון טכנאי")
)

for ($k = 0; $k -lt $pairs.Count; $k++) {
    $left = $pairs[$k][0]
    $right = $pairs[$k][1]

    $leftNode = $null
    $rightNode = $null

    for ($i = 1; $i -le $nodes.Count; $i++) {
        $node = $nodes.Item($i)
        $txt = $node.TextFrame2.TextRange.Text
        if ($txt -eq $left) {
            $leftNode = $node
        } elseif ($txt -eq $right) {
            $rightNode = $node
        }
    }

    if ($leftNode -ne $null) { $leftNode.TextFrame2.TextRange.Text = $right }
    if ($rightNode -ne $null) { $rightNode.TextFrame2.TextRange.Text = $left }
}
